$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was Contact / No display for ContactDetail -> Jurisdiction / United States of America
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 (duplicate Contact / No display for ContactDetail) is removed entirely,
# shifting everything below up by one row.
$meta.Rows(11).Delete()

# Elements sheet: root Extension row (row 2) Short/Definition columns updated
$elements.Range("K2").Value = "Employee Gender"
$elements.Range("L2").Value = "Employee gender code"
